$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-CellText 'D2' '72.267.19'
Set-CellText 'E2' '  +4.31%  '
Set-CellText 'D3' '3.639.48'
Set-CellText 'E3' '  +7.30%  '
Set-CellText 'D4' '1.00'
Set-CellText 'E4' '  +0.04%  '
Set-CellText 'D5' '594.65'
Set-CellText 'E5' '  +1.39%  '
Set-CellText 'D6' '181.45'
Set-CellText 'E6' '  +0.59%  '
Set-CellText 'D7' '3.626.34'
Set-CellText 'E7' '  +7.25%  '
Set-CellText 'D8' '0.608'
Set-CellText 'E8' '  +1.97%  '
Set-CellText 'E9' '  +0.16%  '
Set-CellText 'D10' '0.205'
Set-CellText 'E10' '  +3.99%  '
Set-CellText 'D11' '0.606'
Set-CellText 'E11' '  +2.19%  '
Set-CellText 'D12' '50.11'
Set-CellText 'E12' '  +3.14%  '
Set-CellText 'D13' '0.0000288'
Set-CellText 'E13' '  +1.92%  '
Set-CellText 'D14' '702.87'
Set-CellText 'E14' '  +3.18%  '
Set-CellText 'D15' '4.227.13'
Set-CellText 'E15' '  +7.48%  '
Set-CellText 'D16' '8.95'
Set-CellText 'E16' '  +3.26%  '
Set-CellText 'D17' '72.350.72'
Set-CellText 'E17' '  +4.30%  '
Set-CellText 'D18' '3.638.92'
Set-CellText 'E18' '  +7.36%  '
Set-CellText 'E19' '  +2.28%  '
Set-CellText 'D20' '18.48'
Set-CellText 'E20' '  +4.22%  '
Set-CellText 'D21' '11.64'
Set-CellText 'E21' '  +3.00%  '
Set-CellText 'D22' '0.934'
Set-CellText 'E22' '  +2.91%  '
Set-CellText 'D23' '5.85'
Set-CellText 'E23' '  +8.22%  '
Set-CellText 'D24' '17.86'
Set-CellText 'E24' '  +3.80%  '
Set-CellText 'D25' '104.17'
Set-CellText 'E25' '  +1.03%  '
Set-CellText 'D26' '4.04'
Set-CellText 'E26' '  +2.89%  '
Set-CellText 'D27' '2.86'
Set-CellText 'E27' '  +4.58%  '
Set-CellText 'D28' '9.94'
Set-CellText 'E28' '  +2.64%  '
Set-CellText 'D29' '35.17'
Set-CellText 'E29' '  +3.70%  '
Set-CellText 'D30' '9.11'
Set-CellText 'E30' '  +3.63%  '
Set-CellText 'D31' '7.42'
Set-CellText 'E31' '  +6.64%  '
Set-CellText 'D32' '4.17'
Set-CellText 'E32' '  +16.30%  '
Set-CellText 'D33' '586.55'
Set-CellText 'E33' '  +5.29%  '
Set-CellText 'D34' '11.28'
Set-CellText 'E34' '  +1.15%  '
Set-CellText 'D35' '0.108'
Set-CellText 'E35' '  +1.64%  '
Set-CellText 'D36' '59.84'
Set-CellText 'E36' '  +2.23%  '
Set-CellText 'E37' '  -0.04%  '
Set-CellText 'B38' 'Maker'
Set-CellText 'C38' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText 'D38' '3.636.09'
Set-CellText 'E38' '  -0.98%  '
Set-CellText 'B39' 'Kaspa'
Set-CellText 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-CellText 'D39' '0.144'
Set-CellText 'E39' '  +2.96%  '
Set-CellText 'D40' '0.0₃0779'
Set-CellText 'E40' '  +8.31%  '
Set-CellText 'D41' '35.79'
Set-CellText 'E41' '  +0.25%  '
Set-CellText 'D42' '3.45'
Set-CellText 'E42' '  +5.56%  '
Set-CellText 'D43' '2.78'
Set-CellText 'E43' '  +3.49%  '
Set-CellText 'D44' '0.0455'
Set-CellText 'E44' '  +6.99%  '
Set-CellText 'D45' '0.350'
Set-CellText 'E45' '  +3.28%  '
Set-CellText 'E46' '  +3.93%  '
Set-CellText 'E47' '  +6.05%  '
Set-CellText 'D48' '1.47'
Set-CellText 'E48' '  +5.88%  '
Set-CellText 'D49' '0.132'
Set-CellText 'E49' '  +1.85%  '
Set-CellText 'E50' '  -0.20%  '
Set-CellText 'D51' '133.76'
Set-CellText 'E51' '  -0.25%  '
